# Updated cryptos list on Sat Nov  9 20:31:22 UTC 2024 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $text) {
    # Force the cell to store a literal text value (matches the source
    # inline-string cells) even when the text looks numeric, e.g. "1.00",
    # without leaving the cells number format/style changed afterwards.
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "76.343.07"
$ws.Range("E2").Value = "  -0.23%  "
Set-TextValue $ws.Range("D3") "3.085.84"
$ws.Range("E3").Value = "  +5.00%  "
$ws.Range("E4").Value = "  -0.04%  "
Set-TextValue $ws.Range("D5") "198.31"
$ws.Range("E5").Value = "  -0.11%  "
Set-TextValue $ws.Range("D6") "615.59"
$ws.Range("E6").Value = "  +3.52%  "
Set-TextValue $ws.Range("D7") "0.999"
$ws.Range("E8").Value = "  +0.26%  "
$ws.Range("E9").Value = "  +7.32%  "
Set-TextValue $ws.Range("D10") "3.082.85"
$ws.Range("E10").Value = "  +4.71%  "
$ws.Range("E11").Value = "  -0.01%  "
$ws.Range("E12").Value = "  +0.03%  "
Set-TextValue $ws.Range("D13") "5.23"
$ws.Range("E13").Value = "  +7.49%  "
Set-TextValue $ws.Range("D14") "3.654.22"
$ws.Range("E14").Value = "  +4.88%  "
Set-TextValue $ws.Range("D15") "29.23"
$ws.Range("E15").Value = "  +3.30%  "
Set-TextValue $ws.Range("D16") "76.460.08"
$ws.Range("E16").Value = "  -0.12%  "
$ws.Range("E17").Value = "  +3.45%  "
Set-TextValue $ws.Range("D18") "3.082.24"
$ws.Range("E18").Value = "  +4.26%  "
$ws.Range("E19").Value = "  +0.60%  "
Set-TextValue $ws.Range("D20") "9.15"
$ws.Range("E20").Value = "  +5.64%  "
$ws.Range("B21").Value = "SuiNetwork"
$ws.Range("C21").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
Set-TextValue $ws.Range("D21") "2.58"
$ws.Range("E21").Value = "  +14.43%  "
$ws.Range("B22").Value = "BitcoinCash"
$ws.Range("C22").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
Set-TextValue $ws.Range("D22") "381.46"
$ws.Range("E22").Value = "  +2.16%  "
Set-TextValue $ws.Range("D23") "4.43"
$ws.Range("E23").Value = "  +2.83%  "
Set-TextValue $ws.Range("D24") "3.236.36"
$ws.Range("E24").Value = "  +4.13%  "
Set-TextValue $ws.Range("D25") "72.30"
$ws.Range("E25").Value = "  +0.05%  "
$ws.Range("E26").Value = "  +0.70%  "
Set-TextValue $ws.Range("D27") "4.39"
$ws.Range("E27").Value = "  +3.08%  "
Set-TextValue $ws.Range("D28") "10.02"
$ws.Range("E28").Value = "  +4.39%  "
Set-TextValue $ws.Range("D29") "0.0000109"
$ws.Range("E29").Value = "  +1.83%  "
$ws.Range("E30").Value = "  -0.21%  "
$ws.Range("E31").Value = "  +0.92%  "
$ws.Range("E32").Value = "  +3.58%  "
Set-TextValue $ws.Range("D33") "500.06"
$ws.Range("E33").Value = "  +0.46%  "
Set-TextValue $ws.Range("D34") "1.92"
$ws.Range("E34").Value = "  +5.25%  "
Set-TextValue $ws.Range("D35") "1.00"
$ws.Range("E35").Value = "  -0.05%  "
Set-TextValue $ws.Range("D36") "20.77"
$ws.Range("E36").Value = "  +3.22%  "
$ws.Range("E37").Value = "  +11.90%  "
Set-TextValue $ws.Range("D38") "162.07"
$ws.Range("E38").Value = "  -1.61%  "
Set-TextValue $ws.Range("D39") "195.52"
$ws.Range("E39").Value = "  +8.59%  "
$ws.Range("E40").Value = "  +0.64%  "
Set-TextValue $ws.Range("D41") "0.379"
$ws.Range("E41").Value = "  -3.49%  "
$ws.Range("E42").Value = "  -8.15%  "
Set-TextValue $ws.Range("D44") "0.801"
$ws.Range("E44").Value = "  +22.28%  "
Set-TextValue $ws.Range("D45") "5.13"
$ws.Range("E45").Value = "  +4.85%  "
$ws.Range("E46").Value = "  +6.09%  "
Set-TextValue $ws.Range("D47") "41.33"
$ws.Range("E47").Value = "  +2.97%  "
Set-TextValue $ws.Range("D48") "1.65"
$ws.Range("E48").Value = "  +0.59%  "
Set-TextValue $ws.Range("D49") "2.44"
$ws.Range("E49").Value = "  +6.85%  "
Set-TextValue $ws.Range("D50") "0.600"
$ws.Range("E50").Value = "  +2.11%  "
Set-TextValue $ws.Range("D51") "3.90"
$ws.Range("E51").Value = "  +0.71%  "
